$d = $word.ActiveDocument

# --- Change 1: merge the three "Hidden Markov Model" runs into one run ---
# Text currently reads: " – was used with " + "Hidden Markov Model" + " (HMM)"
# After edit it should read as a single run: " – was used with Hidden Markov Model (HMM)"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    " " + [char]8211 + " was used with Hidden Markov Model (HMM)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " " + [char]8211 + " was used with Hidden Markov Model (HMM)",
    2)

# --- Change 2: add strikethrough to the "K-Means Clustering" paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "K-Means Clustering*") {
        $p.Range.Font.StrikeThrough = 1
        break
    }
}
